$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must remain text
# (to preserve exact formatting, e.g. trailing zeros / precision).
$textForceCells = @(
    @{Cell="D5"; Value="562.94"},
    @{Cell="D6"; Value="143.28"},
    @{Cell="D13"; Value="26.16"},
    @{Cell="D18"; Value="11.34"},
    @{Cell="D19"; Value="323.27"},
    @{Cell="D20"; Value="6.82"},
    @{Cell="D23"; Value="66.75"},
    @{Cell="D24"; Value="1.74"},
    @{Cell="D25"; Value="8.77"},
    @{Cell="D26"; Value="554.14"},
    @{Cell="D28"; Value="0.999"},
    @{Cell="D33"; Value="1.89"},
    @{Cell="D38"; Value="153.21"},
    @{Cell="D39"; Value="5.43"},
    @{Cell="D40"; Value="18.54"},
    @{Cell="D42"; Value="0.993"},
    @{Cell="D43"; Value="147.32"},
    @{Cell="D44"; Value="2.23"},
    @{Cell="D46"; Value="0.0529"},
    @{Cell="D47"; Value="19.82"},
    @{Cell="D48"; Value="0.593"},
    @{Cell="D49"; Value="0.0920"},
    @{Cell="D51"; Value="11.57"}
)

foreach ($item in $textForceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}

# Plain text / non-ambiguous values
$plainCells = @(
    @{Cell="D2"; Value="62.006.05"},
    @{Cell="E2"; Value="  -0.34%  "},
    @{Cell="D3"; Value="2.420.31"},
    @{Cell="E3"; Value="  -0.12%  "},
    @{Cell="E4"; Value="  +0.01%  "},
    @{Cell="E5"; Value="  +0.07%  "},
    @{Cell="E6"; Value="  -0.70%  "},
    @{Cell="E7"; Value="  -0.03%  "},
    @{Cell="E8"; Value="  -0.33%  "},
    @{Cell="E9"; Value="  -0.29%  "},
    @{Cell="E10"; Value="  -0.72%  "},
    @{Cell="E11"; Value="  -4.11%  "},
    @{Cell="E12"; Value="  -0.96%  "},
    @{Cell="E13"; Value="  +0.13%  "},
    @{Cell="E14"; Value="  -1.54%  "},
    @{Cell="D16"; Value="61.945.60"},
    @{Cell="E16"; Value="  -0.10%  "},
    @{Cell="D17"; Value="2.408.69"},
    @{Cell="E17"; Value="  -0.51%  "},
    @{Cell="E18"; Value="  +1.04%  "},
    @{Cell="E19"; Value="  -0.47%  "},
    @{Cell="E20"; Value="  +0.72%  "},
    @{Cell="E21"; Value="  -1.26%  "},
    @{Cell="E22"; Value="  -0.06%  "},
    @{Cell="E23"; Value="  +1.97%  "},
    @{Cell="E24"; Value="  +0.68%  "},
    @{Cell="E25"; Value="  -2.51%  "},
    @{Cell="E26"; Value="  -5.38%  "},
    @{Cell="D27"; Value="2.539.03"},
    @{Cell="E27"; Value="  -0.16%  "},
    @{Cell="E28"; Value="  -2.49%  "},
    @{Cell="D29"; Value="0.0₃0935"},
    @{Cell="E29"; Value="  -0.90%  "},
    @{Cell="E30"; Value="  -1.06%  "},
    @{Cell="E31"; Value="  -4.62%  "},
    @{Cell="E32"; Value="  -1.62%  "},
    @{Cell="E33"; Value="  -0.48%  "},
    @{Cell="E34"; Value="  -3.72%  "},
    @{Cell="E35"; Value="  -0.07%  "},
    @{Cell="E36"; Value="  -0.91%  "},
    @{Cell="E37"; Value="  -1.69%  "},
    @{Cell="E38"; Value="  -1.03%  "},
    @{Cell="E39"; Value="  -4.98%  "},
    @{Cell="E40"; Value="  -0.87%  "},
    @{Cell="E41"; Value="  +0.07%  "},
    @{Cell="E42"; Value="  -0.66%  "},
    @{Cell="E43"; Value="  -2.29%  "},
    @{Cell="E44"; Value="  -6.33%  "},
    @{Cell="E45"; Value="  -0.31%  "},
    @{Cell="E46"; Value="  -2.06%  "},
    @{Cell="B47"; Value="InjectiveProtocol"},
    @{Cell="C47"; Value="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"},
    @{Cell="E47"; Value="  -2.73%  "},
    @{Cell="B48"; Value="Mantle"},
    @{Cell="C48"; Value="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"},
    @{Cell="E48"; Value="  +0.00%  "},
    @{Cell="E49"; Value="  -0.51%  "},
    @{Cell="E50"; Value="  -0.73%  "},
    @{Cell="B51"; Value="WhiteBITCoin"},
    @{Cell="C51"; Value="https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"},
    @{Cell="E51"; Value="  +0.63%  "}
)

foreach ($item in $plainCells) {
    $ws.Range($item.Cell).Value = $item.Value
}
